# TEXAS_2024.xlsx style cleanup:
#  1. Rename header row 1 (A1:D1) to the new machine-friendly column names.
#  2. Title-case the Spanish connector words ("de", "del", "el", "la", "los",
#     "las", "y") inside state/municipality names in columns A and B, except
#     when such a word is the first token of the string. Also normalize the
#     literal "TOTAL" label to "Total".
#  3. Drop the trailing footnote rows (sample size / source / author / date)
#     that followed the final TOTAL row, shrinking the used range back down
#     to A1:D2073.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames ------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Connector-word title-casing -----------------------------------------
# NB: this engine's `-eq`/`-ceq`/`-contains`/`-ccontains` operators are all
# case-INsensitive for strings (no ordinal variant), so exact-case checks
# below go through the case-sensitive .NET `.Equals()` instance method
# instead.
$connectors = @("de", "del", "el", "la", "los", "las", "y")

function Contains-Exact($arr, $val) {
    foreach ($item in $arr) {
        if ($item.Equals($val)) {
            return $true
        }
    }
    return $false
}

function Transform-Label($s) {
    if ($s -eq $null) {
        return $s
    }
    if ($s.Equals("TOTAL")) {
        return "Total"
    }
    $words = $s -split ' '
    for ($i = 0; $i -lt $words.Length; $i++) {
        if ($i -gt 0 -and (Contains-Exact $connectors $words[$i])) {
            $words[$i] = $words[$i].Substring(0, 1).ToUpper() + $words[$i].Substring(1)
        }
    }
    return ($words -join ' ')
}

$lastRow = 2073
for ($r = 2; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($a -ne $null -and $a.Length -gt 0) {
        $newA = Transform-Label($a)
        if (-not $newA.Equals($a)) {
            $ws.Cells.Item($r, 1).Value = $newA
        }
    }
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -ne $null -and $b.Length -gt 0) {
        $newB = Transform-Label($b)
        if (-not $newB.Equals($b)) {
            $ws.Cells.Item($r, 2).Value = $newB
        }
    }
}

# --- 3. Drop trailing footnote rows -----------------------------------------
$ws.Rows("2074:2079").Delete()
